$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (Price and Volume(1h) columns updated by data refresh)
$updates = @{
    'D2' = '246.76'
    'E2' = '1.11%'
    'D3' = '30.10'
    'E3' = '11.12%'
    'D4' = '5.177'
    'E4' = '0.34%'
    'D5' = '0.05739'
    'E5' = '2.07%'
    'D6' = '6.603'
    'E6' = '2.00%'
    'D7' = '0.8615'
    'E7' = '5.55%'
    'D8' = '0.8793'
    'E8' = '5.76%'
    'D9' = '0.1363'
    'E9' = '2.74%'
    'D10' = '0.07080'
    'E10' = '2.57%'
    'D11' = '0.02888'
    'E11' = '0.10%'
    'D12' = '0.09404'
    'E12' = '0.26%'
    'D13' = '0.001518'
    'E13' = '0.61%'
    'D14' = '0.04134'
    'E14' = '-3.32%'
    'D15' = '0.0006006'
    'D16' = '0.005947'
    'E16' = '-2.64%'
    'E17' = '-2.86%'
    'D18' = '3.059'
    'E18' = '1.29%'
    'D19' = '2.269'
    'E19' = '-1.66%'
    'D20' = '0.3185'
    'E20' = '2.31%'
    'D21' = '0.03279'
    'E21' = '5.70%'
    'E22' = '1.28%'
    'D23' = '3.607'
    'E23' = '-3.52%'
    'E24' = '0.42%'
    'D25' = '0.001212'
    'E25' = '-1.04%'
    'D26' = '0.004502'
    'E26' = '0.39%'
    'D27' = '0.0001210'
    'E27' = '23.51%'
    'D28' = '0.0001385'
    'E28' = '0.20%'
    'D40' = '0.03782'
    'E40' = '3.64%'
    'D41' = '0.005708'
    'E41' = '-5.12%'
    'E42' = '1.80%'
    'D43' = '0.002199'
    'E43' = '-15.54%'
    'D44' = '0.01004'
    'E44' = '22.96%'
    'D45' = '0.00005123'
    'E45' = '-3.63%'
    'D46' = '0.00000000750'
    'E46' = '0.03%'
    'D47' = '0.08894'
    'E47' = '-18.34%'
    'D48' = '0.002767'
    'E48' = '5.47%'
    'D49' = '0.00002099'
    'E49' = '0.03%'
    'D50' = '0.0001999'
    'E50' = '0.03%'
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}
